# Update countries & provincias Spain
# Refresh COVID-19 stats for a handful of countries and re-sort a few rows
# that jumped past their neighbours by total-case count (column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 68 / 69 : Costa Rica overtakes Irlanda -----------------------
$ws.Range("A68").Value = "Costa Rica"
$ws.Range("B68").Value = 27737
$ws.Range("C68").Value = 806
$ws.Range("D68").Value = 9010
$ws.Range("E68").Value = 18436
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 10
$ws.Range("H68").Value = 291

$ws.Range("A69").Value = "Irlanda"
$ws.Range("B69").Value = 27191
$ws.Range("C69").Value = 196
$ws.Range("D69").Value = 23364
$ws.Range("E69").Value = 2053
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 1774

# --- Row 92 / 93 / 94 : Libano overtakes Guinea and Gabon --------------
$ws.Range("A92").Value = "Libano"
$ws.Range("B92").Value = 8442
$ws.Range("C92").Value = 397
$ws.Range("D92").Value = 2650
$ws.Range("E92").Value = 5695
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 3
$ws.Range("H92").Value = 97

$ws.Range("A93").Value = "Guinea"
$ws.Range("B93").Value = 8260
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 7177
$ws.Range("E93").Value = 1033
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 50

$ws.Range("A94").Value = "Gabon"
$ws.Range("B94").Value = 8225
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 6277
$ws.Range("E94").Value = 1897
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 51

# --- Row 137 / 138 : Angola overtakes Yemen -----------------------------
$ws.Range("A137").Value = "Angola"
$ws.Range("B137").Value = 1879
$ws.Range("C137").Value = 27
$ws.Range("D137").Value = 628
$ws.Range("E137").Value = 1165
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 86

$ws.Range("A138").Value = "Yemen"
$ws.Range("B138").Value = 1858
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 1009
$ws.Range("E138").Value = 321
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 528

# --- Remaining rows that simply got fresh numbers (no reordering) ------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5513506
$ws.Range("C4").Value = 37240
$ws.Range("D4").Value = 2882936
$ws.Range("E4").Value = 2458269
$ws.Range("G4").Value = 766
$ws.Range("H4").Value = 172301

# Row 6 - India
$ws.Range("B6").Value = 2589208
$ws.Range("C6").Value = 63986
$ws.Range("D6").Value = 1860672
$ws.Range("E6").Value = 678452

# Row 22 - Alemania
$ws.Range("B22").Value = 224446
$ws.Range("C22").Value = 672
$ws.Range("E22").Value = 12606

# Row 116 - Suazilandia
$ws.Range("B116").Value = 3745
$ws.Range("C116").Value = 75
$ws.Range("D116").Value = 2268
$ws.Range("E116").Value = 1408
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 69

# Row 122 - Sri Lanka
$ws.Range("B122").Value = 2890
$ws.Range("C122").Value = 4
$ws.Range("E122").Value = 213

# Row 126 - Mali
$ws.Range("B126").Value = 2614
$ws.Range("C126").Value = 17
$ws.Range("D126").Value = 1986
$ws.Range("E126").Value = 503

# Row 149 - Liberia
$ws.Range("B149").Value = 1257
$ws.Range("C149").Value = 5
$ws.Range("D149").Value = 788
$ws.Range("E149").Value = 387

# Row 158 - Republica del Chad
$ws.Range("B158").Value = 952
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 864
$ws.Range("E158").Value = 12

# --- Timestamp banner ----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 22:11"
